# ============================================================================
# Adds a new worksheet "SPEC - 1.7" (cloned from the TEMPLATE sheet) right
# after "SPEC - 1.6" and before "ORIGINAL SURVEY", documenting the update of
# the Customer class to add an `additionalInformation` field.
# ============================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Clone the TEMPLATE worksheet and put the copy right after "SPEC - 1.6"
# ---------------------------------------------------------------------------
$template = $wb.Worksheets.Item("TEMPLATE")
$spec16   = $wb.Worksheets.Item("SPEC - 1.6")
$template.Copy($null, $spec16)

$ws = $wb.Worksheets.Item("TEMPLATE (2)")
$ws.Name = "SPEC - 1.7"

# ---------------------------------------------------------------------------
# 2. Trim the cloned template rows down to the rows that are actually needed.
#    TEMPLATE layout:
#      1  title bar
#      2  subtitle
#      3  column headers
#      4  first content row   (kept, content replaced)
#      5  second content row  (kept, content replaced)
#      6-13  empty placeholder rows (removed)
#      14 "git add/commit" instructions row
#      15 "git status" result row
#    After deleting rows 6-13, rows 14/15 become rows 6/7.
# ---------------------------------------------------------------------------
$ws.Range("A6:C13").EntireRow.Delete()

# ---------------------------------------------------------------------------
# 3. Selection / active tab bookkeeping
# ---------------------------------------------------------------------------
$ws.Range("A1:C1").Select()

# ---------------------------------------------------------------------------
# 4. Subtitle (row 2)
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "UPDATE THE CUSTOMER CLASS TO HAVE A FIELD - ADDITIONAL INFORMATION"

# ---------------------------------------------------------------------------
# 5. Row 4 - update the constructor
# ---------------------------------------------------------------------------
$code1 = @'
    /* *Constructor
     * @param    {String} customerName          Name of the customer
     * @param    {String} phoneNumber           10 Digit Phone Number of the customer
     * @param    {String} address               150 characters address
     * @param    {String} additionalInformation 250 characters additional information
     * */
    constructor(customerName, phoneNumber, address, additionalInformation) {
        this._customerName = customerName;
        this._phoneNumber = phoneNumber;
        this._address = address;
        this._additionalInformation = additionalInformation;
    }
'@

$stepA   = "Update the constructor of the Customer class to have the additional information field."
$resultA = "Updated the constructor of the Customer class to have the additional information field."

$ws.Range("B4").Value = $code1
$ws.Range("A4").Value = $stepA
$ws.Range("C4").Value = $resultA

# bold the "constructor" / "Customer" keywords within the A4 / C4 sentences
foreach ($addr in @("A4", "C4")) {
    $rng = $ws.Range($addr)
    $txt = $rng.Value
    $i1 = $txt.IndexOf("constructor") + 1
    $rng.Characters($i1, "constructor".Length).Font.Bold = $true
    $i2 = $txt.IndexOf("Customer") + 1
    $rng.Characters($i2, "Customer".Length).Font.Bold = $true
}

# ---------------------------------------------------------------------------
# 6. Row 5 - implement getter / setter
# ---------------------------------------------------------------------------
$code2 = @'
    /* Get the additional information
     * @return  {String} Additional information
     */
    get additionalInformation() {
        return this._additionalInformation;
    }
    /* Set the additional information
     * @param   {String} additionalInformation  250 characters additional information
     */
    set additionalInformation(additionalInformation) {
        this._additionalInformation = additionalInformation;
    }
'@

$stepB   = "Implement the getter and setter methods for the additionalInformation field"
$resultB = "Implemented the getter and setter methods for the additionalInformation field"

$ws.Range("B5").Value = $code2
$ws.Range("A5").Value = $stepB
$ws.Range("C5").Value = $resultB

foreach ($addr in @("A5", "C5")) {
    $rng = $ws.Range($addr)
    $txt = $rng.Value
    $i1 = $txt.IndexOf("additionalInformation") + 1
    $rng.Characters($i1, "additionalInformation".Length).Font.Bold = $true
}

# ---------------------------------------------------------------------------
# 7. Style / layout for rows 4 & 5 (wrap text, left/top or left/center align)
# ---------------------------------------------------------------------------
foreach ($addr in @("A4", "C4", "A5", "C5")) {
    $rng = $ws.Range($addr)
    $rng.HorizontalAlignment = -4131   # xlLeft
    $rng.VerticalAlignment   = -4160   # xlTop
    $rng.WrapText = $true
}
foreach ($addr in @("B4", "B5")) {
    $rng = $ws.Range($addr)
    $rng.HorizontalAlignment = -4131   # xlLeft
    $rng.VerticalAlignment   = -4108   # xlCenter
    $rng.WrapText = $true
    $rng.Font.Color = 6299648   # RGB(0, 32, 96)
}

$ws.Rows.Item(4).RowHeight = 238
$ws.Rows.Item(5).RowHeight = 221

# ---------------------------------------------------------------------------
# 8. Row 6 / 7 - commit instructions & result (re-use template's git-add /
#    git-commit boilerplate wording, only the commit message & git status
#    output are specific to this change)
# ---------------------------------------------------------------------------
$commitMsg = 'git commit -m "Updated the Customer class to have the additionalInformation field."'
$ws.Range("B7").Value = $commitMsg

$gitStatusTail = " - would display the following files.                                                                                                                                                                                                                                          modified:   docs/specification/hafele-v1/HAFELE-SPEC-V-01.xlsx" + "`n        modified:   js/model/customer.js"
$ws.Range("C6").Value = "git status" + $gitStatusTail
$rngC6 = $ws.Range("C6")
$rngC6.Characters(1, "git status".Length).Font.Bold = $true

$ws.Rows.Item(6).RowHeight = 68

# ---------------------------------------------------------------------------
# 9. Sheet protection (matches the other "SPEC - 1.x" sheets)
# ---------------------------------------------------------------------------
$ws.Protect()
